$wb = $excel.ActiveWorkbook

# --- LogInPage: change LogIn action from "Click" to "Enter" for rows 2,3,5,6 ---
$wsLogin = $wb.Worksheets.Item("LogInPage")
$wsLogin.Range("F2").Value = "Enter"
$wsLogin.Range("F3").Value = "Enter"
$wsLogin.Range("F5").Value = "Enter"
$wsLogin.Range("F6").Value = "Enter"

# Update the selection on LogInPage to a single cell (F4) without making it the active sheet
$wsLogin.Range("F4").Select()

# --- Make IndexPage the active sheet (was NewSubmissionPage) ---
$wsIndex = $wb.Worksheets.Item("IndexPage")
$wsIndex.Activate()
